$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.513.77"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "1.659.96"
$ws.Range("E3").Value = "  -3.95%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.510"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.263"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "1.892.33"
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D13").Value = "1.657.13"
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "27.492.59"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "240.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "0.0₃0729"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0499"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").Value = "1.455.77"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.51%  "
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.923"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.573"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.48%  "
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.791"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").Value = "1.801.10"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("E49").Value = "  -6.34%  "
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  -4.43%  "
